$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the first data row) which shifts all subsequent rows up by one.
$ws.Rows.Item(3).Delete()

# Update the selected cell to C13 to match the post-edit selection.
$ws.Range("C13").Select()
